# Sync attendance_reports: normalize "Recorded By" (column G) ordering so
# that "System" is always listed first among the recorded-by entries.
#
# Known transformations observed in the source data:
#   "dnasr281@gmail.com, System"            -> "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com"    -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value = "System, system, backup@backdoor.com"
    }
}
